# Refresh the "cryptos" price table with the latest scraped values.
# Note: several Price (column D) values look like plain numbers (e.g. "1.00",
# "72.70"), but the sheet stores them as literal text so the trailing zeros /
# thousand-dot grouping survive. A bare `.Value = '1.00'` would let Excel's
# usual smart type-detection coerce that into the number 1, so those specific
# assignments are prefixed with a leading apostrophe to force text entry,
# exactly as typing '1.00 into the cell in the Excel UI would.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '42.122.96'
$ws.Range('E2').Value = '  +0.57%  '

# Row 3
$ws.Range('D3').Value = '2.267.81'
$ws.Range('E3').Value = '  +0.06%  '

# Row 4
$ws.Range('B4').Value = 'BinanceUSD'
$ws.Range('C4').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D4').Value = '''218.41'
$ws.Range('E4').Value = '  +21,717.43%  '

# Row 5
$ws.Range('B5').Value = 'TetherUSD'
$ws.Range('C5').Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range('D5').Value = '''1.00'
$ws.Range('E5').Value = '  -0.01%  '

# Row 6
$ws.Range('B6').Value = 'BNB'
$ws.Range('C6').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D6').Value = '''305.69'
$ws.Range('E6').Value = '  +1.23%  '

# Row 7
$ws.Range('B7').Value = 'Solana'
$ws.Range('C7').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D7').Value = '''93.28'
$ws.Range('E7').Value = '  +1.38%  '

# Row 8
$ws.Range('E8').Value = '  -0.35%  '

# Row 9
$ws.Range('E9').Value = '  -0.05%  '

# Row 10
$ws.Range('E10').Value = '  +1.16%  '

# Row 11
$ws.Range('D11').Value = '''32.89'
$ws.Range('E11').Value = '  +2.03%  '

# Row 12
$ws.Range('D12').Value = '''0.0805'
$ws.Range('E12').Value = '  +0.67%  '

# Row 13
$ws.Range('D13').Value = '''0.113'
$ws.Range('E13').Value = '  -1.77%  '

# Row 14
$ws.Range('E14').Value = '  +0.33%  '

# Row 15
$ws.Range('D15').Value = '2.618.87'
$ws.Range('E15').Value = '  +0.06%  '

# Row 16
$ws.Range('D16').Value = '''14.35'
$ws.Range('E16').Value = '  +1.56%  '

# Row 17
$ws.Range('D17').Value = '2.277.41'
$ws.Range('E17').Value = '  +0.88%  '

# Row 18
$ws.Range('E18').Value = '  +3.98%  '

# Row 19
$ws.Range('D19').Value = '41.983.60'
$ws.Range('E19').Value = '  +0.42%  '

# Row 20
$ws.Range('D20').Value = '''12.71'
$ws.Range('E20').Value = '  +5.05%  '

# Row 21
$ws.Range('D21').Value = '0.0₃0920'

# Row 22
$ws.Range('D22').Value = '''5.98'
$ws.Range('E22').Value = '  +0.89%  '

# Row 23
$ws.Range('D23').Value = '''68.21'
$ws.Range('E23').Value = '  +1.92%  '

# Row 24
$ws.Range('D24').Value = '''243.95'
$ws.Range('E24').Value = '  +1.39%  '

# Row 25
$ws.Range('D25').Value = '''2.60'
$ws.Range('E25').Value = '  +1.57%  '

# Row 26
$ws.Range('E26').Value = '  +2.37%  '

# Row 27
$ws.Range('E27').Value = '  -0.15%  '

# Row 28
$ws.Range('D28').Value = '''24.01'
$ws.Range('E28').Value = '  +0.46%  '

# Row 29
$ws.Range('E29').Value = '  +0.61%  '

# Row 30
$ws.Range('D30').Value = '''2.10'
$ws.Range('E30').Value = '  -3.54%  '

# Row 31
$ws.Range('D31').Value = '''35.19'
$ws.Range('E31').Value = '  +4.07%  '

# Row 32
$ws.Range('D32').Value = '''159.87'
$ws.Range('E32').Value = '  +0.51%  '

# Row 33
$ws.Range('E33').Value = '  +3.42%  '

# Row 34
$ws.Range('D34').Value = '''0.999'
$ws.Range('E34').Value = '  -0.02%  '

# Row 35
$ws.Range('E35').Value = '  -0.05%  '

# Row 36
$ws.Range('D36').Value = '''3.05'
$ws.Range('E36').Value = '  -0.81%  '

# Row 37
$ws.Range('D37').Value = '''17.13'
$ws.Range('E37').Value = '  +3.82%  '

# Row 38
$ws.Range('E38').Value = '  -1.03%  '

# Row 39
$ws.Range('E39').Value = '  +1.51%  '

# Row 40
$ws.Range('E40').Value = '  +0.61%  '

# Row 41
$ws.Range('E41').Value = '  -0.15%  '

# Row 42
$ws.Range('D42').Value = '''4.06'
$ws.Range('E42').Value = '  +3.46%  '

# Row 43
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '2.013.60'
$ws.Range('E43').Value = '  -2.37%  '

# Row 44
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = '''19.67'
$ws.Range('E44').Value = '  +0.14%  '

# Row 45
$ws.Range('E45').Value = '  +10.51%  '

# Row 46
$ws.Range('E46').Value = '  +1.67%  '

# Row 47
$ws.Range('E47').Value = '  +1.73%  '

# Row 48
$ws.Range('E48').Value = '  +0.96%  '

# Row 49
$ws.Range('D49').Value = '''53.75'
$ws.Range('E49').Value = '  +4.23%  '

# Row 50
$ws.Range('D50').Value = '''72.70'
$ws.Range('E50').Value = '  +2.89%  '

# Row 51
$ws.Range('E51').Value = '  +0.31%  '
